$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Rights")

$ws.Range("E3").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E4").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_VIEW'
$ws.Range("E5").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E6").Value = 'TASK_DELETE, PERSON_DELETE, ADDITIONAL_TEST_DELETE, SAMPLE_DELETE, IMMUNIZATION_VIEW, CLINICAL_VISIT_DELETE, IMMUNIZATION_DELETE, DOCUMENT_DELETE, PERSON_VIEW, THERAPY_VIEW, TREATMENT_DELETE, PRESCRIPTION_DELETE, SAMPLE_VIEW, PATHOGEN_TEST_DELETE, CASE_VIEW, ADDITIONAL_TEST_VIEW, VISIT_DELETE, CLINICAL_COURSE_VIEW, DOCUMENT_VIEW, TASK_VIEW'
$ws.Range("E7").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E8").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E9").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_EDIT, CASE_VIEW'
$ws.Range("E10").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_EDIT, CASE_VIEW'
$ws.Range("E11").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_EDIT, CASE_VIEW'
$ws.Range("E12").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_EDIT, CASE_VIEW'
$ws.Range("E13").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_EDIT, CASE_VIEW'
$ws.Range("E14").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_EDIT, CASE_VIEW'
$ws.Range("E15").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_EDIT, CASE_VIEW'
$ws.Range("E16").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E17").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_EDIT, CASE_VIEW'
$ws.Range("E19").Value = 'IMMUNIZATION_VIEW, PERSON_VIEW'
$ws.Range("E20").Value = 'IMMUNIZATION_VIEW, PERSON_VIEW, PERSON_EDIT'
$ws.Range("E21").Value = 'IMMUNIZATION_VIEW, PERSON_VIEW'
$ws.Range("E22").Value = 'PERSON_DELETE, IMMUNIZATION_VIEW, VISIT_DELETE, PERSON_VIEW'
$ws.Range("E25").Value = 'VISIT_DELETE, PERSON_VIEW'
$ws.Range("E27").Value = 'PERSON_VIEW, PERSON_EDIT'
$ws.Range("E32").Value = 'ADDITIONAL_TEST_DELETE, SAMPLE_VIEW, PATHOGEN_TEST_DELETE, ADDITIONAL_TEST_VIEW'
$ws.Range("E34").Value = 'SAMPLE_EDIT, SAMPLE_VIEW'
$ws.Range("E35").Value = 'SAMPLE_EDIT, SAMPLE_VIEW'
$ws.Range("E36").Value = 'SAMPLE_EDIT, SAMPLE_VIEW'
$ws.Range("E38").Value = 'SAMPLE_EDIT, SAMPLE_VIEW'
$ws.Range("E44").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E45").Value = 'CONTACT_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E46").Value = 'CONTACT_VIEW, PERSON_VIEW, PERSON_EDIT, CASE_VIEW'
$ws.Range("E47").Value = 'CONTACT_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E48").Value = 'TASK_DELETE, PERSON_DELETE, ADDITIONAL_TEST_DELETE, SAMPLE_DELETE, DOCUMENT_DELETE, PERSON_VIEW, SAMPLE_VIEW, PATHOGEN_TEST_DELETE, ADDITIONAL_TEST_VIEW, CASE_VIEW, CONTACT_VIEW, VISIT_DELETE, DOCUMENT_VIEW, TASK_VIEW'
$ws.Range("E49").Value = 'CONTACT_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E50").Value = 'CONTACT_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E51").Value = 'CONTACT_VIEW, CONTACT_EDIT, PERSON_VIEW, CASE_CREATE, PERSON_EDIT, CASE_VIEW'
$ws.Range("E52").Value = 'CONTACT_VIEW, CONTACT_EDIT, PERSON_VIEW, PERSON_EDIT, CASE_VIEW'
$ws.Range("E53").Value = 'CONTACT_VIEW, CONTACT_EDIT, PERSON_VIEW, PERSON_EDIT, CASE_VIEW'
$ws.Range("E54").Value = 'CONTACT_VIEW, CONTACT_EDIT, PERSON_VIEW, PERSON_EDIT, CASE_VIEW'
$ws.Range("E64").Value = 'TASK_EDIT, TASK_VIEW'
$ws.Range("E66").Value = 'DOCUMENT_DELETE, DOCUMENT_VIEW, EVENT_VIEW'
$ws.Range("E72").Value = 'TASK_DELETE, PERSON_DELETE, ADDITIONAL_TEST_DELETE, SAMPLE_DELETE, DOCUMENT_DELETE, EVENTPARTICIPANT_DELETE, PERSON_VIEW, ACTION_DELETE, SAMPLE_VIEW, PATHOGEN_TEST_DELETE, EVENT_VIEW, ADDITIONAL_TEST_VIEW, VISIT_DELETE, DOCUMENT_VIEW, TASK_VIEW, EVENTPARTICIPANT_VIEW'
$ws.Range("E77").Value = 'PERSON_VIEW, EVENT_VIEW'
$ws.Range("E78").Value = 'PERSON_VIEW, EVENTPARTICIPANT_VIEW, EVENT_VIEW'
$ws.Range("E79").Value = 'PERSON_VIEW, PERSON_EDIT, EVENTPARTICIPANT_VIEW, EVENT_VIEW'
$ws.Range("E80").Value = 'PERSON_VIEW, EVENTPARTICIPANT_VIEW, EVENT_VIEW'
$ws.Range("E81").Value = 'PERSON_DELETE, ADDITIONAL_TEST_DELETE, SAMPLE_DELETE, VISIT_DELETE, PERSON_VIEW, EVENTPARTICIPANT_VIEW, SAMPLE_VIEW, PATHOGEN_TEST_DELETE, ADDITIONAL_TEST_VIEW, EVENT_VIEW'
$ws.Range("E82").Value = 'PERSON_VIEW, EVENTPARTICIPANT_VIEW, EVENT_VIEW'
$ws.Range("E83").Value = 'PERSON_VIEW, PERSON_EDIT, EVENTPARTICIPANT_EDIT, EVENTPARTICIPANT_VIEW, EVENT_VIEW'
$ws.Range("E104").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E105").Value = 'CONTACT_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E106").Value = 'CONTACT_VIEW, PERSON_VIEW, DASHBOARD_CONTACT_VIEW, CASE_VIEW'
$ws.Range("E108").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E109").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E110").Value = 'THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E111").Value = 'THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E112").Value = 'THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E113").Value = 'THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E114").Value = 'THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E115").Value = 'THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E116").Value = 'THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E117").Value = 'CLINICAL_COURSE_VIEW, THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E118").Value = 'CLINICAL_COURSE_VIEW, THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E119").Value = 'CLINICAL_COURSE_VIEW, THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E120").Value = 'CLINICAL_COURSE_VIEW, THERAPY_VIEW, PERSON_VIEW, CASE_VIEW'
$ws.Range("E121").Value = 'PERSON_VIEW, CASE_VIEW'
$ws.Range("E122").Value = 'PERSON_VIEW, PORT_HEALTH_INFO_VIEW, CASE_VIEW'
$ws.Range("E135").Value = 'CAMPAIGN_FORM_DATA_DELETE, CAMPAIGN_VIEW, CAMPAIGN_FORM_DATA_VIEW'
$ws.Range("E142").Value = 'PERSON_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS'
$ws.Range("E143").Value = 'PERSON_VIEW, TRAVEL_ENTRY_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS'
$ws.Range("E144").Value = 'PERSON_VIEW, TRAVEL_ENTRY_VIEW, PERSON_EDIT, TRAVEL_ENTRY_MANAGEMENT_ACCESS'
$ws.Range("E145").Value = 'PERSON_VIEW, TRAVEL_ENTRY_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS'
$ws.Range("E146").Value = 'TASK_DELETE, PERSON_DELETE, DOCUMENT_DELETE, VISIT_DELETE, PERSON_VIEW, TRAVEL_ENTRY_VIEW, DOCUMENT_VIEW, TASK_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS'
$ws.Range("E165").Value = 'PERSON_DELETE, IMMUNIZATION_VIEW, IMMUNIZATION_DELETE, IMMUNIZATION_EDIT, PERSON_VIEW, CASE_CREATE, EVENTPARTICIPANT_EDIT, SAMPLE_CREATE, PATHOGEN_TEST_CREATE, VISIT_DELETE, EVENT_EDIT, IMMUNIZATION_CREATE, PERSON_EDIT, PATHOGEN_TEST_EDIT, CASE_EDIT, EXTERNAL_MESSAGE_VIEW, CONTACT_CREATE, PATHOGEN_TEST_DELETE, SAMPLE_VIEW, CASE_VIEW, EVENT_VIEW, SAMPLE_EDIT, CONTACT_VIEW, CONTACT_EDIT, EVENT_CREATE, EVENTPARTICIPANT_VIEW, EVENTPARTICIPANT_CREATE'

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = '1.0.0'

